$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$rowData = @{
    2 = @(1.02, 1.007280975910814, 1.032557869945718, 1.009980696467818, 1.019235045214043, 1.031824640371513, 1.012553531049677, 1.035362755253459, 1.012852271110154, 1.022078895016406, 1.00828416421713)
    3 = @(1.02, 1.009231370332541, 1.033190643017001, 1.011657959203368, 1.021653614723384, 1.032103544889807, 1.014129963555196, 1.035804796420873, 1.014331060234558, 1.024299017094116, 1.008809098177585)
    4 = @(1.02, 1.010483495590752, 1.033593845226922, 1.012735060703889, 1.023191958683203, 1.032276666823106, 1.015140085664442, 1.036083910111618, 1.015279319716399, 1.025708751043483, 1.009145302381481)
    5 = @(1.02, 1.011007560670793, 1.033761865532169, 1.01318594758285, 1.023832383872346, 1.032347696831334, 1.015562400652198, 1.036199602500027, 1.015675939839855, 1.026295055448888, 1.009285826336596)
    6 = @(1.02, 1.011095418224029, 1.03378998999512, 1.013261541533608, 1.023939547253839, 1.032359520674142, 1.015633173041526, 1.036218931423452, 1.015742416152548, 1.026393128454254, 1.009309373423025)
    7 = @(1.02, 1.010490507264614, 1.033596096149165, 1.012741092998537, 1.023200540703013, 1.032277622796263, 1.015145737803643, 1.036085462458968, 1.015284627297653, 1.025716610114044, 1.009147183260029)
    8 = @(1.02, 1.007942199917082, 1.032773017121833, 1.010549255383179, 1.020057976543921, 1.031920424327065, 1.013088374342414, 1.035513582771571, 1.013353839986636, 1.022834798766358, 1.008462293030155)
    9 = @(1.02, 1.00337379917891, 1.031274382094596, 1.006622461135669, 1.014312574460263, 1.031234293847877, 1.009385164455902, 1.034452452900573, 1.009883947749325, 1.017547561299422, 1.007228308149759)
    10 = @(1.02, 1.000272803540536, 1.030242168666808, 1.003958771814581, 1.010336962073828, 1.030738103098054, 1.006861409425307, 1.03370846622473, 1.007522896789279, 1.013876829237041, 1.006386549387825)
    11 = @(1.02, 0.998916204734638, 1.02978718622272, 1.002793920926661, 1.00857967168641, 1.030513891376141, 1.005754956427701, 1.033377475559993, 1.006488653839517, 1.012251495084205, 1.006017325018511)
    12 = @(1.02, 0.9984101612681098, 1.029616964015273, 1.002359471083953, 1.00792143350884, 1.030429188375601, 1.005341865518827, 1.033253187668991, 1.006102654443775, 1.011642268864958, 1.005879448815885)
    13 = @(1.02, 0.998518807271824, 1.029653532789806, 1.00245274309825, 1.008062878813373, 1.030447421999275, 1.005430571071211, 1.033279908881378, 1.006185536491193, 1.011773200944679, 1.005909057054694)
    14 = @(1.02, 0.9988744190922028, 1.029773140612609, 1.002758045626582, 1.008525374473934, 1.030506918872696, 1.005720853467221, 1.033367229376763, 1.006456784742672, 1.012201249228042, 1.006005943124385)
    15 = @(1.02, 0.9990932375086428, 1.029846672603907, 1.002945916168736, 1.008809600519427, 1.030543388131654, 1.005899425440834, 1.03342085195574, 1.00662366502538, 1.012464250975716, 1.006065540565665)
    16 = @(1.02, 1.000362539562343, 1.030272193930344, 1.004035833441335, 1.010452822751194, 1.030752784884015, 1.006934548893071, 1.03373024547212, 1.007591281406489, 1.013983931713174, 1.006410952254793)
    17 = @(1.02, 1.00115498896606, 1.030536952482703, 1.004716407142918, 1.011473898474637, 1.030881618362493, 1.007580162315451, 1.033921942838366, 1.008195023823701, 1.014927501793225, 1.006626339018693)
    18 = @(1.02, 1.001615879541936, 1.03069060825557, 1.005112272512519, 1.01206602658349, 1.030955862909642, 1.007955424408571, 1.034032904861146, 1.008546032525908, 1.015474416042096, 1.006751514656973)
    19 = @(1.02, 1.001772807268429, 1.030742870172838, 1.005247067019232, 1.012267345666863, 1.030981025849758, 1.008083157936056, 1.034070596007121, 1.008665524824554, 1.015660317136884, 1.006794119538223)
    20 = @(1.02, 1.001070104834426, 1.030508626502234, 1.004643502349875, 1.011364704215827, 1.030867889135908, 1.007511030330386, 1.033901463738561, 1.008130366656153, 1.014826623744802, 1.00660327733337)
    21 = @(1.02, 0.9987697599199871, 1.029737952937443, 1.002668190988459, 1.008389333911072, 1.030489437870947, 1.00563543109071, 1.033341552873305, 1.00637695996403, 1.012075352574895, 1.005977432887551)
    22 = @(1.02, 0.9973110149106639, 1.029246324471902, 1.001415953192112, 1.006486702454022, 1.030243262286456, 1.004443960108152, 1.032981735023057, 1.005263878275918, 1.010313610194625, 1.005579707571113)
    23 = @(1.02, 0.9980855224426762, 1.029507622107674, 1.002080780307278, 1.007498388522827, 1.030374549772223, 1.005076757190106, 1.033173224012514, 1.005854969801251, 1.011250607151006, 1.005790956672896)
    24 = @(1.02, 1.00110846443161, 1.030521428182789, 1.004676448289491, 1.011414055048986, 1.030874095568423, 1.007542272153907, 1.033910719993691, 1.00815958597955, 1.014872216859775, 1.006613699323814)
    25 = @(1.02, 1.004564378265805, 1.031667600895255, 1.0076455221356, 1.015823086429389, 1.03141845804117, 1.010352021319737, 1.034733174019662, 1.010789243764403, 1.018939749993393, 1.007550623888328)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}

Write-Output "Applied $($rowData.Keys.Count) rows"